$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New benchmark rows (write in the same order the original author entered
# them so new shared-string indices line up with the target workbook)
$ws.Range("A4").Value = 156
$ws.Range("B4").Value = " split per 10 files,extract with read_parquet, write to csv, and manual query to copy from csv"
$ws.Range("C4").Value = "397 s"
$ws.Range("D4").Value = "avg 270 mb"

$ws.Range("A5").Value = 156
$ws.Range("B5").Value = "to_csv loop all file"

$ws.Range("A6").Value = 156
$ws.Range("D6").Value = "avg 350 mb"
$ws.Range("B6").Value = "write manual to csv"
$ws.Range("C6").Value = "376 s"

$ws.Range("D5").Value = "avg 300 mb"
$ws.Range("C5").Value = "294 s"

# Column B now holds a much longer string - widen it to fit
$ws.Columns.Item(2).ColumnWidth = 78

# Scroll/selection state left behind after entering the data
$ws.Range("D21").Select()
